# Weekly price-sheet update: a new Sandia ("Primera", Peru origin) price
# record for a more recent date (2023-10-16, serial 45215) is inserted as
# row 118, pushing the existing rows 118-159 down to 119-160. All other
# rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 118; rows 118:159 shift down to 119:160
# and formatting (e.g. the date cell's number format) is inherited from the
# row above, same as interactively inserting a row in Excel.
$ws.Rows.Item(118).Insert()

$ws.Range("A118").Value = 8
$ws.Range("B118").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C118").Value = 'Coquimbo'
$ws.Range("D118").Value = 45215
$ws.Range("E118").Value = 4
$ws.Range("F118").Value = 100112028
$ws.Range("G118").Value = 'Sandia'
$ws.Range("H118").Value = 'Sin especificar'
$ws.Range("I118").Value = 'Primera'
$ws.Range("J118").Value = 1200
$ws.Range("K118").Value = 750
$ws.Range("L118").Value = 800
$ws.Range("M118").Value = 775
$ws.Range("N118").Value = '$/kilo (volumen en unidades)'
$ws.Range("O118").Value = 'Perú'
$ws.Range("P118").Value = 775
$ws.Range("Q118").Value = 1
$ws.Range("R118").Value = 'Hortaliza'
